$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 9 with the AC Milan vs Pisa match data
$ws.Range("A9").Value = "24/10/2025"
$ws.Range("B9").Value = "AC Milan"
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = "Pisa"
$ws.Range("F9").Value = "D"
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 2
$ws.Range("K9").Value = 1.74
$ws.Range("L9").Value = 1.05
$ws.Range("M9").Value = 20
$ws.Range("N9").Value = 8
$ws.Range("O9").Value = 7
$ws.Range("P9").Value = 2
